$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.869.04'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '3.825.51'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.17'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.48'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '3.825.70'
$ws.Range("E7").Value = '  -2.01%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("E13").Value = '  +8.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.09'
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("D15").Value = '4.467.27'
$ws.Range("E15").Value = '  -2.07%  '
$ws.Range("D16").Value = '3.821.81'
$ws.Range("E16").Value = '  -2.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.68'
$ws.Range("E17").Value = '  +2.79%  '
$ws.Range("D18").Value = '68.000.77'
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.48'
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.88'
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '469.04'
$ws.Range("E22").Value = '  -0.93%  '
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000151'
$ws.Range("E24").Value = '  -9.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.82'
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.26'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.45'
$ws.Range("E28").Value = '  +3.94%  '
$ws.Range("E30").Value = '  -1.58%  '
$ws.Range("D31").Value = '3.972.29'
$ws.Range("E31").Value = '  -2.07%  '
$ws.Range("E32").Value = '  -1.62%  '
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.78'
$ws.Range("E34").Value = '  -2.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.36'
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("D36").Value = '3.791.24'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.90'
$ws.Range("E37").Value = '  +4.75%  '
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.99'
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.02'
$ws.Range("E40").Value = '  -1.69%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.139'
$ws.Range("E41").Value = '  -1.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("E43").Value = '  +1.46%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("B45").Value = 'Cosmos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.81'
$ws.Range("E45").Value = '  +1.57%  '
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '411.11'
$ws.Range("E47").Value = '  -4.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.61'
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000285'
$ws.Range("E49").Value = '  -5.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '143.39'
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("E51").Value = '  -0.16%  '
